$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of data (2-11, excluding 9 which is unchanged) were reshuffled.
# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
$rows = @{
    2  = @{ D = 44369; J = 25; K = 8000;  L = 8000;  M = 8000;  P = 800  }
    3  = @{ D = 44525; J = 20; K = 9000;  L = 9000;  M = 9000;  P = 900  }
    4  = @{ D = 44473; J = 25; K = 11000; L = 11000; M = 11000; P = 1100 }
    5  = @{ D = 44425; J = 30; K = 13000; L = 13000; M = 13000; P = 1300 }
    6  = @{ D = 44469; J = 20; K = 12000; L = 12000; M = 12000; P = 1200 }
    7  = @{ D = 44348; J = 20; K = 10000; L = 10000; M = 10000; P = 1000 }
    8  = @{ D = 44530; J = 30; K = 10000; L = 10000; M = 10000; P = 1000 }
    10 = @{ D = 44526; J = 25; K = 9000;  L = 9000;  M = 9000;  P = 900  }
    11 = @{ D = 44523; J = 30; K = 9000;  L = 9000;  M = 9000;  P = 900  }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals.P  # P: Precio $/Kg
}
